$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4, col B: "Investigate embedding YouTube videos on pages"
#   -> "Put back videos sections with links to YouTube pages if Killa agrees"
$ws.Range("B4").Value = "Put back videos sections with links to YouTube pages if Killa agrees"

# Row 13, col B: "Remove code that's commented out"
#   -> "Remove code that's not used any more"
$ws.Range("B13").Value = "Remove code that's not used any more"

# Row 15, col B: "Featured songs, videos etc"
#   -> "Featured content section(s)"
$ws.Range("B15").Value = "Featured content section(s)"

# Update the selection to match the committed state (B5 single cell selected)
$ws.Range("B5").Select()
